# The two occurrence records on row 2 and row 3 had been mismatched:
# the "Antal"/"Enhet" (count/unit) observation fields, the coordinates,
# the accuracy and the record Id had ended up on the wrong row.
# This swaps those fields back between row 2 and row 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose row-2/row-3 values must be swapped:
#   A  - Id
#   I  - Antal (count)
#   J  - Enhet (unit)
#   Q  - Ost (easting)
#   R  - Nord (northing)
#   S  - Noggrannhet (accuracy)
$columns = @("A", "I", "J", "Q", "R", "S")

foreach ($col in $columns) {
    $rangeRow2 = $ws.Range($col + "2")
    $rangeRow3 = $ws.Range($col + "3")

    $valueRow2 = $rangeRow2.Value2
    $valueRow3 = $rangeRow3.Value2

    $rangeRow2.Value2 = $valueRow3
    $rangeRow3.Value2 = $valueRow2
}
